$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = 'Record'
$ws.Range("B47").Value = 'Balanço Geral'
$ws.Range("C47").Value = 'Social'
$ws.Range("D47").Value = '2025-04-04T13:17'
$ws.Range("E47").Value = 'Positivo'
$ws.Range("F47").Value = 'Termina hoje o mutirão de atualização do Cadastro Único em Campos. Repórter *ao vivo*. Desde quarta-feira, o mutirão é realizado. No primeiro dia, confusão. Muitas pessoas e 800 vagas. Ontem, estava organizado. Distribuíram pulseirinhas separando grupo prioritário. Hoje, seguiu tranquilamente. Fundação de Esportes estava aberta para que as pessoas pudessem se abrigar. Tudo mais tranquilo. Mutirões são realizados para desafogar o atendimento. '

$ws.Range("A48").Value = 'Record'
$ws.Range("B48").Value = 'Balanço Geral'
$ws.Range("C48").Value = 'Obras'
$ws.Range("D48").Value = '2025-04-04T13:21'
$ws.Range("E48").Value = 'Negativo'
$ws.Range("F48").Value = 'Moradores denunciam continuidade de obras irregulares em praça do Pq. Lebret. Repórter *ao vivo* informa que a emissora divulgou a denúncia no dia 26 de março. Secretaria de Obras disse que enviaria equipe ao local. Dez dias depois, situação ainda não foi resolvida. Obra já está em fase final. Repórter questiona como a população pode usar o local, se a praça está com obra sem autorização da prefeitura.  '

$ws.Range("A49").Value = 'Record'
$ws.Range("B49").Value = 'Balanço Geral'
$ws.Range("C49").Value = 'Iluminação'
$ws.Range("D49").Value = '2025-04-04T13:36'
$ws.Range("E49").Value = 'Negativo'
$ws.Range("F49").Value = 'Tô na Bronca. Moradores de Balança Rangel cobram melhorias na iluminação pública. *nota coberta*. Moradores da região passam na estrada sem iluminação, desde o período das eleições. Já solicitaram troca de lâmpadas e, até agora, nada. Exibido vídeo enviado por moradora, mostrando a escuridão na rua que dá acesso à localidade. Não tem acostamento e está cheia de buracos.  Equipe procurou a prefeitura, mas até o momento, nenhuma resposta. '

"done"